# Auto-generated script applying scheduled market-data refresh to Asura_Profits workbook.
# For each changed sheet/row, the price/profit columns (H-N) are overwritten with freshly
# pulled values. A couple of rows gain/lose a trailing column cell (N) to match the refresh.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 17
$ws.Cells.Item(17, 8).Value = 501.05554
$ws.Cells.Item(17, 10).Value = 370.7143
$ws.Cells.Item(17, 12).Value = 1112.1429
$ws.Cells.Item(17, 14).Value = -1448.1429
# Row 70
$ws.Cells.Item(70, 8).Value = 92262
$ws.Cells.Item(70, 9).Value = 201216.4
$ws.Cells.Item(70, 10).Value = 1466.6666
$ws.Cells.Item(70, 11).Value = 603649.2
$ws.Cells.Item(70, 12).Value = 4399.9998
$ws.Cells.Item(70, 13).Value = -603379.2
$ws.Cells.Item(70, 14).Value = -4939.9998
# Row 73
$ws.Cells.Item(73, 8).Value = 92262
$ws.Cells.Item(73, 9).Value = 201216.4
$ws.Cells.Item(73, 10).Value = 1466.6666
$ws.Cells.Item(73, 11).Value = 603649.2
$ws.Cells.Item(73, 12).Value = 4399.9998
$ws.Cells.Item(73, 13).Value = -602713.2
$ws.Cells.Item(73, 14).Value = -6271.9998
# Row 88
$ws.Cells.Item(88, 8).Value = 5000
$ws.Cells.Item(88, 10).Value = 6000
$ws.Cells.Item(88, 12).Value = 6000
$ws.Cells.Item(88, 14).Value = -6812
# Row 91
$ws.Cells.Item(91, 8).Value = 5000
$ws.Cells.Item(91, 10).Value = 6000
$ws.Cells.Item(91, 12).Value = 6000
$ws.Cells.Item(91, 14).Value = -8808
# Row 113
$ws.Cells.Item(113, 8).Value = 3082.1
$ws.Cells.Item(113, 9).Value = 2505
$ws.Cells.Item(113, 10).Value = 3329.4285
$ws.Cells.Item(113, 11).Value = 2505
$ws.Cells.Item(113, 12).Value = 3329.4285
$ws.Cells.Item(113, 13).Value = 749
$ws.Cells.Item(113, 14).Value = -9837.4285
# Row 129
$ws.Cells.Item(129, 8).Value = 1078.1757
$ws.Cells.Item(129, 10).Value = 1126.3478
$ws.Cells.Item(129, 12).Value = 3379.0434
$ws.Cells.Item(129, 14).Value = -13379.0434
# Row 132
$ws.Cells.Item(132, 8).Value = 1644.814
$ws.Cells.Item(132, 9).Value = 1333.1082
$ws.Cells.Item(132, 11).Value = 3999.3246
$ws.Cells.Item(132, 13).Value = -1469.3246
# Row 138
$ws.Cells.Item(138, 8).Value = 3317.5051
$ws.Cells.Item(138, 9).Value = 1973.5
$ws.Cells.Item(138, 10).Value = 3919.2986
$ws.Cells.Item(138, 11).Value = 5920.5
$ws.Cells.Item(138, 12).Value = 11757.8958
$ws.Cells.Item(138, 13).Value = -780.5
$ws.Cells.Item(138, 14).Value = -22037.8958

$ws = $wb.Worksheets.Item("ARM")
# Row 102
$ws.Cells.Item(102, 8).Value = 4588.75
$ws.Cells.Item(102, 9).Value = 4815.7144
$ws.Cells.Item(102, 10).Value = 3000
$ws.Cells.Item(102, 11).Value = 4815.7144
$ws.Cells.Item(102, 12).Value = 3000
$ws.Cells.Item(102, 13).Value = -3193.7144
$ws.Cells.Item(102, 14).Value = -6244
# Row 112
$ws.Cells.Item(112, 8).Value = 20371.297
$ws.Cells.Item(112, 10).Value = 20371.297
$ws.Cells.Item(112, 12).Value = 20371.297
$ws.Cells.Item(112, 14).Value = -23325.297
# Row 122
$ws.Cells.Item(122, 8).Value = 6076.96
$ws.Cells.Item(122, 9).Value = 6121.25
$ws.Cells.Item(122, 10).Value = 5014
$ws.Cells.Item(122, 11).Value = 18363.75
$ws.Cells.Item(122, 12).Value = 15042
$ws.Cells.Item(122, 13).Value = -15913.75
$ws.Cells.Item(122, 14).Value = -19942
# Row 124
$ws.Cells.Item(124, 8).Value = 39444.25
$ws.Cells.Item(124, 10).Value = 39444.25
$ws.Cells.Item(124, 12).Value = 39444.25
$ws.Cells.Item(124, 14).Value = -49264.25
# Row 135
$ws.Cells.Item(135, 8).Value = 26967.727
$ws.Cells.Item(135, 10).Value = 26967.727
$ws.Cells.Item(135, 12).Value = 26967.727
$ws.Cells.Item(135, 14).Value = -37107.727
# Row 125
$ws.Cells.Item(125, 8).Value = 0
$ws.Cells.Item(125, 10).Value = 0
$ws.Cells.Item(125, 12).Value = 0
$ws.Cells.Item(125, 14).Value = ""

$ws = $wb.Worksheets.Item("CRP")
# Row 62
$ws.Cells.Item(62, 8).Value = 34230
$ws.Cells.Item(62, 9).Value = 41476.152
$ws.Cells.Item(62, 11).Value = 41476.152
$ws.Cells.Item(62, 13).Value = -40852.152
# Row 65
$ws.Cells.Item(65, 8).Value = 34230
$ws.Cells.Item(65, 9).Value = 41476.152
$ws.Cells.Item(65, 11).Value = 207380.76
$ws.Cells.Item(65, 13).Value = -204260.76
# Row 122
$ws.Cells.Item(122, 8).Value = 8254.166999999999
$ws.Cells.Item(122, 9).Value = 9702.200000000001
$ws.Cells.Item(122, 10).Value = 1014
$ws.Cells.Item(122, 11).Value = 29106.6
$ws.Cells.Item(122, 12).Value = 3042
$ws.Cells.Item(122, 13).Value = -26656.6
$ws.Cells.Item(122, 14).Value = -7942

$ws = $wb.Worksheets.Item("CUL")
# Row 37
$ws.Cells.Item(37, 8).Value = 54000
$ws.Cells.Item(37, 10).Value = 54000
$ws.Cells.Item(37, 12).Value = 162000
$ws.Cells.Item(37, 14).Value = -162224
# Row 68
$ws.Cells.Item(68, 8).Value = 1170.8392
$ws.Cells.Item(68, 9).Value = 834
$ws.Cells.Item(68, 10).Value = 1423.4688
$ws.Cells.Item(68, 11).Value = 2502
$ws.Cells.Item(68, 12).Value = 4270.4064
$ws.Cells.Item(68, 13).Value = -1691
$ws.Cells.Item(68, 14).Value = -5892.4064
# Row 71
$ws.Cells.Item(71, 8).Value = 1170.8392
$ws.Cells.Item(71, 9).Value = 834
$ws.Cells.Item(71, 10).Value = 1423.4688
$ws.Cells.Item(71, 11).Value = 7506
$ws.Cells.Item(71, 12).Value = 12811.2192
$ws.Cells.Item(71, 13).Value = -3450
$ws.Cells.Item(71, 14).Value = -20923.2192
# Row 103
$ws.Cells.Item(103, 8).Value = 2269
$ws.Cells.Item(103, 10).Value = 2116.6
$ws.Cells.Item(103, 12).Value = 6349.799999999999
$ws.Cells.Item(103, 14).Value = -8107.799999999999
# Row 117
$ws.Cells.Item(117, 8).Value = 85183.336
$ws.Cells.Item(117, 9).Value = 533.3333
$ws.Cells.Item(117, 10).Value = 169833.33
$ws.Cells.Item(117, 11).Value = 1599.9999
$ws.Cells.Item(117, 12).Value = 509499.99
$ws.Cells.Item(117, 13).Value = 1842.0001
$ws.Cells.Item(117, 14).Value = -516383.99
# Row 119
$ws.Cells.Item(119, 8).Value = 6241.1177
$ws.Cells.Item(119, 9).Value = 2788.7778
$ws.Cells.Item(119, 10).Value = 10125
$ws.Cells.Item(119, 11).Value = 8366.3334
$ws.Cells.Item(119, 12).Value = 30375
$ws.Cells.Item(119, 13).Value = -3528.3334
$ws.Cells.Item(119, 14).Value = -40051
# Row 131
$ws.Cells.Item(131, 8).Value = 22879.912
$ws.Cells.Item(131, 10).Value = 26853.514
$ws.Cells.Item(131, 12).Value = 80560.542
$ws.Cells.Item(131, 14).Value = -90640.542

$ws = $wb.Worksheets.Item("GSM")
# Row 102
$ws.Cells.Item(102, 8).Value = 3617.3333
$ws.Cells.Item(102, 9).Value = 3604.2856
$ws.Cells.Item(102, 10).Value = 3800
$ws.Cells.Item(102, 11).Value = 3604.2856
$ws.Cells.Item(102, 12).Value = 3800
$ws.Cells.Item(102, 13).Value = -1982.2856
$ws.Cells.Item(102, 14).Value = -7044
# Row 113
$ws.Cells.Item(113, 8).Value = 1805.7
$ws.Cells.Item(113, 9).Value = 1299.8
$ws.Cells.Item(113, 10).Value = 2311.6
$ws.Cells.Item(113, 11).Value = 1299.8
$ws.Cells.Item(113, 12).Value = 2311.6
$ws.Cells.Item(113, 13).Value = 870.2
$ws.Cells.Item(113, 14).Value = -6651.6

$ws = $wb.Worksheets.Item("LTW")
# Row 22
$ws.Cells.Item(22, 8).Value = 1255.4286
$ws.Cells.Item(22, 9).Value = 745
$ws.Cells.Item(22, 10).Value = 1459.6
$ws.Cells.Item(22, 11).Value = 745
$ws.Cells.Item(22, 12).Value = 1459.6
$ws.Cells.Item(22, 13).Value = -450
$ws.Cells.Item(22, 14).Value = -2049.6
# Row 27
$ws.Cells.Item(27, 8).Value = 1255.4286
$ws.Cells.Item(27, 9).Value = 745
$ws.Cells.Item(27, 10).Value = 1459.6
$ws.Cells.Item(27, 11).Value = 745
$ws.Cells.Item(27, 12).Value = 1459.6
$ws.Cells.Item(27, 13).Value = -638
$ws.Cells.Item(27, 14).Value = -1673.6
# Row 40
$ws.Cells.Item(40, 8).Value = 1796
$ws.Cells.Item(40, 9).Value = 1750
$ws.Cells.Item(40, 11).Value = 1750
$ws.Cells.Item(40, 13).Value = -1614
# Row 61
$ws.Cells.Item(61, 8).Value = 14624.412
$ws.Cells.Item(61, 9).Value = 18185.77
$ws.Cells.Item(61, 10).Value = 3050
$ws.Cells.Item(61, 11).Value = 18185.77
$ws.Cells.Item(61, 12).Value = 3050
$ws.Cells.Item(61, 13).Value = -17983.77
$ws.Cells.Item(61, 14).Value = -3454
# Row 113
$ws.Cells.Item(113, 8).Value = 14624.412
$ws.Cells.Item(113, 9).Value = 18185.77
$ws.Cells.Item(113, 10).Value = 3050
$ws.Cells.Item(113, 11).Value = 18185.77
$ws.Cells.Item(113, 12).Value = 3050
$ws.Cells.Item(113, 13).Value = -16015.77
$ws.Cells.Item(113, 14).Value = -7390
# Row 132
$ws.Cells.Item(132, 8).Value = 3876.5
$ws.Cells.Item(132, 9).Value = 3764.487
$ws.Cells.Item(132, 10).Value = 5332.6665
$ws.Cells.Item(132, 11).Value = 11293.461
$ws.Cells.Item(132, 12).Value = 15997.9995
$ws.Cells.Item(132, 13).Value = -8763.460999999999
$ws.Cells.Item(132, 14).Value = -21057.9995

$ws = $wb.Worksheets.Item("WVR")
# Row 122
$ws.Cells.Item(122, 8).Value = 41669668
$ws.Cells.Item(122, 9).Value = 50003000
$ws.Cells.Item(122, 10).Value = 3005
$ws.Cells.Item(122, 11).Value = 150009000
$ws.Cells.Item(122, 12).Value = 9015
$ws.Cells.Item(122, 13).Value = -150006550
$ws.Cells.Item(122, 14).Value = -13915
